$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns starting at column S. This shifts the existing
# S:AE columns (sat_act_required.* through act scores.75th) three
# positions to the right, becoming V:AH, and grows the used range from
# A1:AE2 to A1:AH2.
$ws.Range("S1:U1").EntireColumn.Insert()

# Column R used to hold "general_college_subjects.arts"; it is relabeled
# as "history" while the three freshly inserted columns pick up the new
# "electives"/"cs" subjects plus the relocated "arts" header.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"
$ws.Range("U1").Value = "general_college_subjects.arts"

# Populate the data row for the three newly inserted subject columns.
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0

# Normalize the "Unknown" placeholder text to lowercase "unknown" for the
# importance columns.
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"
